$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 370.69446
$ws.Range("J33").Value = 455.25
$ws.Range("L33").Value = 455.25
$ws.Range("N33").Value = -913.25
$ws.Range("H51").Value = 2434.7144
$ws.Range("I51").Value = 1980.6666
$ws.Range("J51").Value = 2775.25
$ws.Range("K51").Value = 1980.6666
$ws.Range("L51").Value = 2775.25
$ws.Range("M51").Value = -1496.6666
$ws.Range("N51").Value = -3743.25
$ws.Range("H101").Value = 1000
$ws.Range("I101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("M101").ClearContents()
$ws.Range("H105").Value = 19900
$ws.Range("J105").Value = 19900
$ws.Range("L105").Value = 19900
$ws.Range("N105").Value = -26888
$ws.Range("H107").Value = 3048.1365
$ws.Range("I107").Value = 3521.25
$ws.Range("J107").Value = 2480.4
$ws.Range("K107").Value = 3521.25
$ws.Range("L107").Value = 2480.4
$ws.Range("M107").Value = -1601.25
$ws.Range("N107").Value = -6320.4
$ws.Range("H132").Value = 9016620
$ws.Range("I132").Value = 16675467
$ws.Range("K132").Value = 50026401
$ws.Range("M132").Value = -50023871

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 933.0769
$ws.Range("I2").Value = 768.3333
$ws.Range("J2").Value = 1074.2858
$ws.Range("K2").Value = 768.3333
$ws.Range("L2").Value = 1074.2858
$ws.Range("M2").Value = -655.3333
$ws.Range("N2").Value = -1300.2858
$ws.Range("H45").Value = 1054.1538
$ws.Range("I45").Value = 964
$ws.Range("K45").Value = 964
$ws.Range("M45").Value = -587
$ws.Range("H63").Value = 1945.0625
$ws.Range("I63").Value = 1788.2727
$ws.Range("K63").Value = 1788.2727
$ws.Range("M63").Value = -1102.2727
$ws.Range("H66").Value = 1945.0625
$ws.Range("I66").Value = 1788.2727
$ws.Range("K66").Value = 8941.363499999999
$ws.Range("M66").Value = -5509.363499999999
$ws.Range("H116").Value = 933.0769
$ws.Range("I116").Value = 768.3333
$ws.Range("J116").Value = 1074.2858
$ws.Range("K116").Value = 768.3333
$ws.Range("L116").Value = 1074.2858
$ws.Range("M116").Value = 1525.6667
$ws.Range("N116").Value = -5662.2858

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 933.0769
$ws.Range("I3").Value = 768.3333
$ws.Range("J3").Value = 1074.2858
$ws.Range("K3").Value = 768.3333
$ws.Range("L3").Value = 1074.2858
$ws.Range("M3").Value = -654.3333
$ws.Range("N3").Value = -1302.2858
$ws.Range("H107").Value = 1297.6428
$ws.Range("I107").Value = 993
$ws.Range("K107").Value = 993
$ws.Range("M107").Value = 927

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1307.4
$ws.Range("J31").Value = 1391.6666
$ws.Range("L31").Value = 1391.6666
$ws.Range("N31").Value = -1981.6666
$ws.Range("H34").Value = 1307.4
$ws.Range("J34").Value = 1391.6666
$ws.Range("L34").Value = 1391.6666
$ws.Range("N34").Value = -1795.6666
$ws.Range("H86").Value = 3051104
$ws.Range("I86").Value = 5565738.5
$ws.Range("K86").Value = 5565738.5
$ws.Range("M86").Value = -5564615.5
$ws.Range("H89").Value = 3051104
$ws.Range("I89").Value = 5565738.5
$ws.Range("K89").Value = 27828692.5
$ws.Range("M89").Value = -27823076.5
$ws.Range("H99").Value = 1445.6428
$ws.Range("I99").Value = 1431.1428
$ws.Range("J99").Value = 1460.1428
$ws.Range("K99").Value = 1431.1428
$ws.Range("L99").Value = 1460.1428
$ws.Range("M99").Value = 66.85719999999992
$ws.Range("N99").Value = -4456.1428
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H126").Value = 1445.6428
$ws.Range("I126").Value = 1431.1428
$ws.Range("J126").Value = 1460.1428
$ws.Range("K126").Value = 4293.428400000001
$ws.Range("L126").Value = 4380.428400000001
$ws.Range("M126").Value = -1823.428400000001
$ws.Range("N126").Value = -9320.428400000001
$ws.Range("H132").Value = 2095.875
$ws.Range("I132").Value = 1775.9166
$ws.Range("K132").Value = 5327.7498
$ws.Range("M132").Value = -2797.7498
$ws.Range("H134").Value = 2026
$ws.Range("I134").Value = 1984.174
$ws.Range("K134").Value = 5952.522
$ws.Range("M134").Value = -3417.522

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 3016.1428
$ws.Range("J87").Value = 3819.8
$ws.Range("L87").Value = 11459.4
$ws.Range("N87").Value = -13955.4
$ws.Range("H88").Value = 5904.95
$ws.Range("I88").Value = 2000
$ws.Range("J88").Value = 6110.4736
$ws.Range("K88").Value = 6000
$ws.Range("L88").Value = 18331.4208
$ws.Range("M88").Value = -5572
$ws.Range("N88").Value = -19187.4208
$ws.Range("H90").Value = 3016.1428
$ws.Range("J90").Value = 3819.8
$ws.Range("L90").Value = 34378.2
$ws.Range("N90").Value = -46858.2
$ws.Range("H91").Value = 5904.95
$ws.Range("I91").Value = 2000
$ws.Range("J91").Value = 6110.4736
$ws.Range("K91").Value = 6000
$ws.Range("L91").Value = 18331.4208
$ws.Range("M91").Value = -4518
$ws.Range("N91").Value = -21295.4208
$ws.Range("H96").Value = 7685.7144
$ws.Range("J96").Value = 8215.385
$ws.Range("L96").Value = 24646.155
$ws.Range("N96").Value = -28764.155
$ws.Range("H102").Value = 4099.75
$ws.Range("I102").Value = 1999
$ws.Range("J102").Value = 4800
$ws.Range("K102").Value = 5997
$ws.Range("L102").Value = 14400
$ws.Range("M102").Value = -3563
$ws.Range("N102").Value = -19268
$ws.Range("H131").Value = 11496602
$ws.Range("J131").Value = 2620.169
$ws.Range("L131").Value = 7860.507
$ws.Range("N131").Value = -17940.507

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 844.125
$ws.Range("I102").Value = 807.5714
$ws.Range("J102").Value = 1100
$ws.Range("K102").Value = 807.5714
$ws.Range("L102").Value = 1100
$ws.Range("M102").Value = 814.4286
$ws.Range("N102").Value = -4344
$ws.Range("H113").Value = 1171.3334
$ws.Range("I113").Value = 1311.8334
$ws.Range("J113").Value = 1101.0834
$ws.Range("K113").Value = 1311.8334
$ws.Range("L113").Value = 1101.0834
$ws.Range("M113").Value = 858.1666
$ws.Range("N113").Value = -5441.0834
$ws.Range("H126").Value = 2221.3
$ws.Range("I126").Value = 1833.1666
$ws.Range("K126").Value = 5499.4998
$ws.Range("M126").Value = -3029.4998

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2021.25
$ws.Range("I7").Value = 1695
$ws.Range("K7").Value = 1695
$ws.Range("M7").Value = -1583
$ws.Range("H16").Value = 680.05884
$ws.Range("I16").Value = 749.4286
$ws.Range("K16").Value = 749.4286
$ws.Range("M16").Value = -579.4286
$ws.Range("H40").Value = 2224.1
$ws.Range("I40").Value = 2109.889
$ws.Range("J40").Value = 3252
$ws.Range("K40").Value = 2109.889
$ws.Range("L40").Value = 3252
$ws.Range("M40").Value = -1973.889
$ws.Range("N40").Value = -3524
$ws.Range("H126").Value = 2021.25
$ws.Range("I126").Value = 1695
$ws.Range("K126").Value = 5085
$ws.Range("M126").Value = -2615

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 521.3125
$ws.Range("I107").Value = 421.9091
$ws.Range("K107").Value = 1265.7273
$ws.Range("M107").Value = 654.2727
$ws.Range("H126").Value = 500012500
$ws.Range("I126").Value = 1000000000
$ws.Range("J126").Value = 25000
$ws.Range("K126").Value = 3000000000
$ws.Range("L126").Value = 75000
$ws.Range("M126").Value = -2999997530
$ws.Range("N126").Value = -79940
$ws.Range("H136").Value = 867.125
$ws.Range("I136").Value = 758
$ws.Range("J136").Value = 1085.375
$ws.Range("K136").Value = 2274
$ws.Range("L136").Value = 3256.125
$ws.Range("M136").Value = 276
$ws.Range("N136").Value = -8356.125
